{"js": "// Fix the typo \"a opportunity\" -> \"an opportunity\" in the apology letter's\n// closing paragraph (\"... and grant me a opportunity to take final exam...\").\nconst body = context.document.body;\n\nconst results = body.search(\"a opportunity\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Target text \"a opportunity\" not found in document body.');\n}\n\n// Replace the matched text in place; the engine keeps this run's existing\n// character formatting (Times New Roman, 12pt) for the replacement text.\nresults.items[0].insertText(\"an opportunity\", \"Replace\");\nawait context.sync();\n", "ps1": "# Fix the typo \"a opportunity\" -> \"an opportunity\" in the apology letter's\n# closing paragraph (\"... and grant me a opportunity to take final exam...\").\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"a opportunity\"\n$find.Replacement.Text = \"an opportunity\"\n$find.Forward = $true\n$find.Wrap = 0\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n\n$found = $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $false, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2)\n\nif (-not $found) {\n    throw 'Target text \"a opportunity\" not found in document content.'\n}\n"}
